$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New benchmark rows to append below the existing TGCN row (A2:N2):
# columns -> A:Model B:seq_len C:pre_len D:learning_rate E:hidden_dim
#            F:MAE G:MSE H:RMSE I:R2 J:Accuracy K:val_loss L:epochs N:How-to-run note
$rows = @(
    @{ Row=3;  A="TGCN"; B=32; C=7; D=0.0001; E=64; F=840.8;   G=3376470;  H=1837.5; I=0.95883828401565496; J=0.8168; K=298655514624; L=100;  N=$null },
    @{ Row=4;  A="TGCN"; B=32; C=7; D=0.005;  E=64; F=694.12;  G=2294786;  H=1514;   I=0.97;               J=0.848;  K=202978459648; L=100;  N=$null },
    @{ Row=5;  A="TGCN"; B=32; C=7; D=0.001;  E=64; F=613.23;  G=1980323;  H=1407;   I=0.975;              J=0.859;  K=175163588608; L=300;  N=$null },
    @{ Row=6;  A="TGCN"; B=32; C=7; D=0.002;  E=64; F=623.3;   G=2049461;  H=1431;   I=0.975;              J=0.857;  K=181278982144; L=200;  N=$null },
    @{ Row=7;  A="TGCN"; B=32; C=7; D=0.1;    E=64; F=725.98;  G=2251792;  H=1500;   I=0.972;              J=0.85;   K=199175536640; L=100;  N=$null },
    @{ Row=8;  A="GRU";  B=32; C=7; D=0.001;  E=64; F=690.74;  G=3438364;  H=1854;   I=0.958;              J=0.815;  K=3438364;      L=100;  N="Folder: TGCN/T-GCN/T-GCN-PyTorch:  python main.py --data shenzhen --model_name GRU --max_epochs 100 --learning_rate 0.0001 --weight_decay 0 --batch_size 32 --hidden_dim 64 --loss mse_with_regularizer --settings supervised --gpus 0 --seq_len 32 --pre_len 7" },
    @{ Row=9;  A="GRU";  B=32; C=7; D=0.005;  E=64; F=625.05;  G=2184623;  H=1478;   I=0.973;              J=0.852;  K=2184623;      L=100;  N=$null },
    @{ Row=10; A="GRU";  B=32; C=7; D=0.01;   E=64; F=606.38;  G=2482061;  H=1575;   I=0.969;              J=0.842;  K=2482061;      L=100;  N=$null },
    @{ Row=11; A="GCN";  B=32; C=7; D=0.001;  E=64; F=1969;    G=22413590; H=4734;   I=0.726;              J=0.527;  K=22413590;     L=100;  N=$null },
    @{ Row=12; A="GCN";  B=32; C=7; D=0.001;  E=64; F=1942;    G=22150196; H=4706;   I=0.73;               J=0.53;   K=22150196;     L=3000; N=$null }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    if ($r.N -ne $null) {
        $ws.Cells.Item($rowNum, 14).Value = $r.N
    }
}

$ws.Range("W9").Select()
